$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value must be forced
# to Text (the source sheet stores these numeric-looking prices/percentages as
# plain text, so a bare assignment would otherwise be auto-coerced to a Number).
$updates = @(
    @{ Cell = 'D2'; Value = '44.481.16'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +0.56%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.239.05'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -0.32%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  +0.48%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '303.93'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -1.05%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '94.52'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -1.84%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  -0.80%  '; ForceText = $false },
    @{ Cell = 'E8'; Value = '  +0.22%  '; ForceText = $false },
    @{ Cell = 'E9'; Value = '  -2.32%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '34.89'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -1.95%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '7.14'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  -1.62%  '; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -0.13%  '; ForceText = $false },
    @{ Cell = 'B14'; Value = 'WrappedEther'; ForceText = $false },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; ForceText = $false },
    @{ Cell = 'D14'; Value = '2.375.73'; ForceText = $false },
    @{ Cell = 'E14'; Value = '  +1.77%  '; ForceText = $false },
    @{ Cell = 'B15'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false },
    @{ Cell = 'D15'; Value = '2.580.23'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -0.25%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '0.829'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -0.93%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '13.50'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -1.04%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '44.300.80'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  +0.52%  '; ForceText = $false },
    @{ Cell = 'E19'; Value = '  -3.37%  '; ForceText = $false },
    @{ Cell = 'E20'; Value = '  -3.80%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '11.70'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -4.18%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '65.07'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -0.80%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '237.17'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -0.03%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '2.93'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -1.37%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '1.97'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -1.88%  '; ForceText = $false },
    @{ Cell = 'E26'; Value = '  +0.04%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +5.03%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '9.68'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -3.47%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '37.17'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -2.52%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '19.83'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -1.74%  '; ForceText = $false },
    @{ Cell = 'E31'; Value = '  -2.53%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '149.92'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -1.72%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '0.0781'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -2.43%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  +0.26%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  -2.23%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '1.87'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  +5.87%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  -2.12%  '; ForceText = $false },
    @{ Cell = 'E38'; Value = '  -0.92%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '14.82'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +1.04%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '3.35'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -5.85%  '; ForceText = $false },
    @{ Cell = 'E41'; Value = '  -2.81%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '0.0295'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -1.26%  '; ForceText = $false },
    @{ Cell = 'E43'; Value = '  +0.21%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '1.819.07'; ForceText = $false },
    @{ Cell = 'E44'; Value = '  +3.27%  '; ForceText = $false },
    @{ Cell = 'E45'; Value = '  +11.16%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '79.92'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -3.88%  '; ForceText = $false },
    @{ Cell = 'E47'; Value = '  -2.79%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '97.93'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -2.55%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '4.82'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -2.83%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '68.36'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +0.46%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '7.93'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -3.01%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Mimic typing a leading-apostrophe / "format as text" entry in the UI so
        # Excel keeps the literal digits instead of parsing them into a Number.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
